$d = $word.ActiveDocument

# 1. Expand the intro bullet for the ConstructConnect director role.
$d.Content.Find.Execute(
    "Responsible for modernizing software delivery practices at ConstructConnect to improve quality. This work ranges from improving how test cases are written to direct, frequent interfacing with executive business leadership to ensure their goals are clear and met.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Responsible for modernizing software delivery practices at ConstructConnect to improve quality. Works with multiple groups including software engineering, operations support, product ownership, and others. Constantly collaborates with executives across the entire organization.",
    2) | Out-Null

# 2. Combine the staff-management sentence into one.
$d.Content.Find.Execute(
    "Manages 19 full-time staff software testers. Oversees 30 additional external testing consultants.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Manages 19 full-time staff software testers and 30 additional external testing consultants.",
    2) | Out-Null

# 3. Fix typo "mulitple" -> "multiple".
$d.Content.Find.Execute(
    "In first six months was directly responsible for reducing post-code freeze regression testing from mulitple weeks to three days by focusing testers on risk-based analysis as well as eliminating duplicate and outdated manual test cases.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In first six months was directly responsible for reducing post-code freeze regression testing from multiple weeks to three days by focusing testers on risk-based analysis as well as eliminating duplicate and outdated manual test cases.",
    2) | Out-Null

# 4. Insert five new organizational-impact bullets before the "Oversaw implementation
#    of pairwise testing..." bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Oversaw implementation of pairwise testing approaches via Hexawise*") {
        $r = $p.Range.Duplicate
        $r.Collapse(1)
        $r.InsertBefore("Mentored and raised up two senior testers into managerial positions.`r")
        $r.InsertBefore("Created a career path ladder for testers with well-defined criteria tied to explicit learning plans. This plan was adopted by other organizational groups within ConstructConnect.`r")
        $r.InsertBefore("Collaborated with Infosec team to bring security scanning into release process, blocking releases with high-level vulnerabilities.`r")
        $r.InsertBefore("Deeply involved as a mentor with ConstructConnect’s internal mentoring program.`r")
        $r.InsertBefore("Works directly with software engineers to raise software craftsmanship and quality capabilities through mentoring, workshops, and building communities of learning.`r")
        break
    }
}

# 5. Remove the "Drove implementation of C4 model architectural diagrams..." bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Drove implementation of C4 model architectural diagrams*") {
        $p.Range.Delete()
        break
    }
}

# 6. Remove the "Driving adoption of static code analysis tooling..." bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Driving adoption of static code analysis tooling*") {
        $p.Range.Delete()
        break
    }
}

# 7. Remove the "Restructured QA department..." bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Restructured QA department*") {
        $p.Range.Delete()
        break
    }
}

# 8. Fix typo "organziations" -> "organizations".
$d.Content.Find.Execute(
    "Led client organziations to improve their software delivery practices. Helped clients improve their realization of driving business value by working with them through executive consulting, coaching, and hands-on software delivery practices.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Led client organizations to improve their software delivery practices. Helped clients improve their realization of driving business value by working with them through executive consulting, coaching, and hands-on software delivery practices.",
    2) | Out-Null
